$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Columns.Item(5).ColumnWidth = 20.17
$ws.Range("K2").Value = 2588964.198201469
$ws.Range("K3").Value = 1275622.430198984
$ws.Range("K4").Value = 1275197.430198984
$ws.Range("K5").Value = 2165280.644397186
$ws.Range("K6").Value = 0.07196287856928146
$ws.Range("K7").Value = 1.651605656346414
$ws.Range("H8").Value = 14.37070569252619
$ws.Range("K8").Value = 0.09131141048528836
$ws.Range("K11").Value = 28353128.97306071
$ws.Range("E12").Value = 0.01841357768863145
$ws.Range("K13").Value = 263910.7235679377
$ws.Range("K14").Value = 514.9341300699099
$ws.Range("K15").Value = 55.58421341009927
$ws.Range("K16").Value = 9.264035568349879
$ws.Range("K17").Value = 0.07632548095692555
$ws.Range("H18").Value = 119.6033070112643
$ws.Range("K18").Value = 1275622.430198984
$ws.Range("H19").Value = 38.09416660875394
$ws.Range("K19").Value = 423683.5538042828
$ws.Range("K20").Value = 466051.9091847111
$ws.Range("H21").Value = 23.19945864779369
$ws.Range("K21").Value = 419584.810607416
$ws.Range("H22").Value = 7.931438853946561
$ws.Range("K22").Value = 4098.743196866766
$ws.Range("H23").Value = 2.379431656183968
$ws.Range("K23").Value = 0.8206740338348483
$ws.Range("K24").Value = 23.73474478423399
$ws.Range("K25").Value = 2588964.198201469
$ws.Range("H27").Value = 5.653692311274727
$ws.Range("H30").Value = 18.68882058728235
$ws.Range("H31").Value = 66.89841397345116
$ws.Range("H32").Value = 0.9342075400657516
$ws.Range("H35").Value = 64.76173353994677
$ws.Range("H36").Value = 37.06600542627907
$ws.Range("H38").Value = 9.856094576956949
$ws.Range("H39").Value = 7.730270256436822
$ws.Range("H40").Value = 5.411189179505775
$ws.Range("E42").Value = 13.23433652621411
$ws.Range("E43").Value = 5.293734610485646
$ws.Range("H43").Value = 6.638937984939859
$ws.Range("E44").Value = 11.9109028735927
$ws.Range("H45").Value = 67.04790083132906
$ws.Range("E46").Value = 29.18048965842116
$ws.Range("H46").Value = 0.9205945069050846
$ws.Range("E47").Value = 9.831221419473341
$ws.Range("E48").Value = 514.9341300699099
$ws.Range("E49").Value = 55.58421341009927
$ws.Range("H49").Value = 32.01664710221132

$ws = $wb.Worksheets.Item(2)
$ws.Range("K2").Value = 2572020.246254603
$ws.Range("K3").Value = 1269638.620066455
$ws.Range("K4").Value = 1269213.620066455
$ws.Range("K5").Value = 2159279.890312709
$ws.Range("K6").Value = 0.111829776455133
$ws.Range("K7").Value = 1.765175244692074
$ws.Range("K8").Value = 0.09132671779534859
$ws.Range("K11").Value = 28162845.5324341
$ws.Range("E12").Value = 0.02030005954817041
$ws.Range("K13").Value = 262183.5113409381
$ws.Range("H14").Value = 14.13731058997421
$ws.Range("K14").Value = 511.7584949430478
$ws.Range("K15").Value = 71.53729761062043
$ws.Range("K16").Value = 7.153729761062041
$ws.Range("K17").Value = 0.07437328679825092
$ws.Range("K18").Value = 1269638.620066455
$ws.Range("H19").Value = 156.2894430915086
$ws.Range("K19").Value = 412740.3559418936
$ws.Range("H20").Value = 22.37266501988687
$ws.Range("K20").Value = 454014.391536083
$ws.Range("K21").Value = 408852.9946258218
$ws.Range("H22").Value = 26.51985094060275
$ws.Range("K22").Value = 3887.361316071823
$ws.Range("H23").Value = 9.066615706188974
$ws.Range("K23").Value = 0.8296150510912133
$ws.Range("H24").Value = 2.719984711856692
$ws.Range("K24").Value = 2572020.246254603
$ws.Range("K25").Value = 24.95485372195421
$ws.Range("H28").Value = 6.46286965723214
$ws.Range("H31").Value = 13.46478894403764
$ws.Range("H32").Value = 42.2340451262468
$ws.Range("H33").Value = 0.840443527069524
$ws.Range("H36").Value = 70.27649321185062
$ws.Range("H37").Value = 21.84482757518156
$ws.Range("H39").Value = 7.259984153487387
$ws.Range("H40").Value = 11.38821043684296
$ws.Range("E41").Value = 10.21961394437434
$ws.Range("H41").Value = 7.971747305790071
$ws.Range("E42").Value = 4.087845577749738
$ws.Range("E43").Value = 15.32942091656152
$ws.Range("E44").Value = 20.48899923060518
$ws.Range("H44").Value = 9.780463081053366
$ws.Range("E45").Value = 19.83202404846683
$ws.Range("E46").Value = 7.591713215820942
$ws.Range("H46").Value = 41.05263028283417
$ws.Range("E47").Value = 511.7584949430478
$ws.Range("H47").Value = 0.8304561089610522
$ws.Range("E48").Value = 71.53729761062043
$ws.Range("H50").Value = 21.93620008524906
$ws.Range("H52").Value = 20.3628356137201

$ws = $wb.Worksheets.Item(3)
$ws.Range("K2").Value = 2576612.94395303
$ws.Range("K3").Value = 1271260.54576063
$ws.Range("K4").Value = 1270835.54576063
$ws.Range("K5").Value = 2160906.408704584
$ws.Range("K6").Value = 0.1557028617907942
$ws.Range("K7").Value = 1.659216450185598
$ws.Range("K8").Value = 0.09537320004849245
$ws.Range("K11").Value = 27016110.84290925
$ws.Range("E12").Value = 0.01873621119097947
$ws.Range("K13").Value = 262651.6762439378
$ws.Range("H14").Value = 14.83276342201457
$ws.Range("K14").Value = 515.6046474057424
$ws.Range("K15").Value = 64.22489532296599
$ws.Range("K16").Value = 8.028111915370749
$ws.Range("K17").Value = 0.07490242024942624
$ws.Range("K18").Value = 1271260.54576063
$ws.Range("H19").Value = 102.2506970670776
$ws.Range("K19").Value = 415706.535248447
$ws.Range("H20").Value = 38.66461917669569
$ws.Range("K20").Value = 457277.1887732917
$ws.Range("K21").Value = 411761.8051058073
$ws.Range("H22").Value = 21.45059758612448
$ws.Range("K22").Value = 3944.730142639659
$ws.Range("H23").Value = 7.333537636281872
$ws.Range("K23").Value = 0.8307664987870182
$ws.Range("H24").Value = 2.200061290884562
$ws.Range("K24").Value = 2576612.94395303
$ws.Range("K25").Value = 24.61076526887319
$ws.Range("H28").Value = 5.227496058682975
$ws.Range("H32").Value = 67.44463345148401
$ws.Range("H33").Value = 0.9382295180405841
$ws.Range("H36").Value = 76.43684226252361
$ws.Range("H37").Value = 36.33361409745591
$ws.Range("H39").Value = 10.70771980366434
$ws.Range("H40").Value = 8.923099836386946
$ws.Range("E41").Value = 11.4687313076725
$ws.Range("H41").Value = 5.353859901832167
$ws.Range("E42").Value = 4.587492523069
$ws.Range("E43").Value = 13.762477569207
$ws.Range("E44").Value = 29.34969811739325
$ws.Range("H44").Value = 7.287198199716006
$ws.Range("E45").Value = 28.61242253332858
$ws.Range("E46").Value = 8.519628971413857
$ws.Range("H46").Value = 65.85507125137894
$ws.Range("E47").Value = 515.6046474057424
$ws.Range("H47").Value = 0.9071818107580653
$ws.Range("E48").Value = 64.22489532296599
$ws.Range("H50").Value = 31.77593709834741

$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 32.22366166980459
